# Generate Report for Handoff
# Update the "Latest Handoff Datetime" (column D) for the row corresponding to the
# e924d709-8e56-4d53-bdb9-cc448585afe6 file (row 5) on both the zh-cn and de-de
# localization-status sheets, recording the new handoff timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-03-03 12:16:06"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-03 12:16:18"
